$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Sep 25, 2024"
$ws.Range("B5").Value = 75100
$ws.Range("C5").Value = 10674.13
$ws.Range("D5").Value = 9446.139999999999
$ws.Range("E5").Value = 7.0246
